$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.1579721434076752
$ws.Range("C2").Value = 0.9336505161243889
$ws.Range("D2").Value = 1.74911684130149
$ws.Range("E2").Value = 1.322541810795216
$ws.Range("F2").Value = 1.328253901032809
$ws.Range("G2").Value = 44

# Row 3 (Q0)
$ws.Range("B3").Value = 0.246962814551667
$ws.Range("C3").Value = 1.238503523070974
$ws.Range("D3").Value = 3.975336168250551
$ws.Range("E3").Value = 1.993824507886928
$ws.Range("F3").Value = 1.985678071014797
$ws.Range("G3").Value = 138

# Row 4 (Q1)
$ws.Range("B4").Value = 0.3026363429728921
$ws.Range("C4").Value = 1.343823136482059
$ws.Range("D4").Value = 8.464273444294301
$ws.Range("E4").Value = 2.90934244190922
$ws.Range("F4").Value = 2.915397665348333
$ws.Range("G4").Value = 67
